$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.336.41"
$ws.Range("E2").Value = "  -0.27%  "

$ws.Range("D3").Value = "2.491.82"
$ws.Range("E3").Value = "  -0.88%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'568.56"
$ws.Range("E5").Value = "  -0.82%  "

$ws.Range("D6").Value = "'166.03"
$ws.Range("E6").Value = "  -0.25%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -0.70%  "

$ws.Range("D9").Value = "'0.159"
$ws.Range("E9").Value = "  -0.21%  "

$ws.Range("E10").Value = "  -0.79%  "

$ws.Range("D11").Value = "'0.348"
$ws.Range("E11").Value = "  -2.58%  "

$ws.Range("D12").Value = "'4.87"
$ws.Range("E12").Value = "  -0.53%  "

$ws.Range("D13").Value = "2.946.11"
$ws.Range("E13").Value = "  -0.93%  "

$ws.Range("D14").Value = "69.222.64"
$ws.Range("E14").Value = "  -0.15%  "

$ws.Range("D15").Value = "'0.0000175"
$ws.Range("E15").Value = "  -0.96%  "

$ws.Range("D16").Value = "'24.09"
$ws.Range("E16").Value = "  -2.87%  "

$ws.Range("D17").Value = "2.484.90"
$ws.Range("E17").Value = "  -1.26%  "

$ws.Range("D18").Value = "'11.18"
$ws.Range("E18").Value = "  -0.71%  "

$ws.Range("D19").Value = "'353.64"
$ws.Range("E19").Value = "  +1.03%  "

$ws.Range("D20").Value = "'7.38"
$ws.Range("E20").Value = "  -2.30%  "

$ws.Range("D21").Value = "'3.90"
$ws.Range("E21").Value = "  -0.18%  "

$ws.Range("E22").Value = "  -4.12%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").Value = "'69.12"
$ws.Range("E24").Value = "  -1.58%  "

$ws.Range("D25").Value = "'3.80"
$ws.Range("E25").Value = "  -3.32%  "

$ws.Range("D26").Value = "2.617.23"
$ws.Range("E26").Value = "  -0.89%  "

$ws.Range("D27").Value = "'8.61"
$ws.Range("E27").Value = "  -3.65%  "

$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("D29").Value = "0.0₃0870"
$ws.Range("E29").Value = "  -2.31%  "

$ws.Range("D30").Value = "'7.52"
$ws.Range("E30").Value = "  -3.79%  "

$ws.Range("D31").Value = "'3.68"
$ws.Range("E31").Value = "  +143.49%  "

$ws.Range("D32").Value = "'438.65"
$ws.Range("E32").Value = "  -5.26%  "

$ws.Range("D33").Value = "'1.19"
$ws.Range("E33").Value = "  -3.42%  "

$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("E35").Value = "  -1.32%  "

$ws.Range("D36").Value = "'154.49"
$ws.Range("E36").Value = "  -1.89%  "

$ws.Range("D37").Value = "'0.113"
$ws.Range("E37").Value = "  -3.30%  "

$ws.Range("D38").Value = "'19.05"
$ws.Range("E38").Value = "  -0.19%  "

$ws.Range("D39").Value = "'18.12"
$ws.Range("E39").Value = "  -1.97%  "

$ws.Range("D41").Value = "'0.313"
$ws.Range("E41").Value = "  -1.75%  "

$ws.Range("D42").Value = "'4.59"
$ws.Range("E42").Value = "  -2.37%  "

$ws.Range("E43").Value = "  -1.93%  "

$ws.Range("D44").Value = "'2.18"
$ws.Range("E44").Value = "  -0.96%  "

$ws.Range("E45").Value = "  -4.45%  "

$ws.Range("D46").Value = "'138.28"
$ws.Range("E46").Value = "  -2.27%  "

$ws.Range("D47").Value = "'3.43"
$ws.Range("E47").Value = "  -1.19%  "

$ws.Range("D48").Value = "'0.504"
$ws.Range("E48").Value = "  -3.11%  "

$ws.Range("D49").Value = "'0.0722"
$ws.Range("E49").Value = "  -1.48%  "

$ws.Range("E50").Value = "  -0.69%  "

$ws.Range("E51").Value = "  -0.56%  "
